# Assign myself with "Enter test result" task.
# Row 4 of Sheet1 is the "Enter Test results" task. Column D is "Assignees",
# column E is "Initial " (the week-4 initial estimate/remaining amount).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = "Md Mostafizur Rahman"
$ws.Range("E4").Value = 5

# Leave the cursor on the cell that was just edited, matching the saved
# selection left behind in the workbook.
[void]$ws.Range("E4").Select()
